# Update the presenter-names textbox ("Textfeld 3") on the last slide (slide 30).
# Original text: "Pascal "
# New text (two runs): "Pascal\t\tSimon \t" + "\tPhilipp "
# The shape is also widened to fit the new, longer text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Textfeld 3") {
        $targetShape = $sh
        break
    }
}
if ($targetShape -eq $null) {
    $targetShape = $s.Shapes.Item(7)
}

$tr = $targetShape.TextFrame.TextRange
$tr.Text = "Pascal`t`tSimon `t"
$run2 = $tr.InsertAfter("`tPhilipp ")

# Widen the textbox to accommodate the added names (EMU -> points: /12700)
$targetShape.Width = 3911648 / 12700
